$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "39.784.62"
$ws.Cells.Item(2, 5).Value = "  -0.84%  "

$ws.Cells.Item(3, 4).Value = "2.226.24"
$ws.Cells.Item(3, 5).Value = "  -5.16%  "

$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "294.36"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -5.29%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "84.70"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.87%  "

$ws.Cells.Item(7, 5).Value = "  -2.71%  "

$ws.Cells.Item(8, 5).Value = "  -0.02%  "

$ws.Cells.Item(9, 5).Value = "  -3.36%  "

$ws.Cells.Item(10, 5).Value = "  -2.97%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "29.95"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.80%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "47.92"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -8.72%  "

$ws.Cells.Item(13, 5).Value = "  -2.46%  "

$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.33"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.51%  "

$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "2.569.23"
$ws.Cells.Item(15, 5).Value = "  -5.14%  "

$ws.Cells.Item(16, 5).Value = "  -4.66%  "

$ws.Cells.Item(17, 4).Value = "2.233.50"
$ws.Cells.Item(17, 5).Value = "  -5.73%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.722"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -5.33%  "

$ws.Cells.Item(19, 4).Value = "39.708.60"
$ws.Cells.Item(19, 5).Value = "  -0.93%  "

$ws.Cells.Item(21, 5).Value = "  -5.40%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "65.30"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -4.36%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.53"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -1.41%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "232.67"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.11%  "

$ws.Cells.Item(25, 5).Value = "  +0.03%  "

$ws.Cells.Item(26, 5).Value = "  -5.56%  "

$ws.Cells.Item(27, 5).Value = "  +0.65%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "22.85"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.98%  "

$ws.Cells.Item(29, 5).Value = "  -0.74%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.21"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.99%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "32.55"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -6.40%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "151.57"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.50%  "

$ws.Cells.Item(33, 5).Value = "  -0.20%  "

$ws.Cells.Item(34, 5).Value = "  -5.90%  "

$ws.Cells.Item(35, 5).Value = "  -2.22%  "

$ws.Cells.Item(36, 5).Value = "  -4.26%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "16.03"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +2.53%  "

$ws.Cells.Item(38, 5).Value = "  -1.95%  "

$ws.Cells.Item(39, 5).Value = "  -0.83%  "

$ws.Cells.Item(40, 5).Value = "  -5.28%  "

$ws.Cells.Item(41, 5).Value = "  -4.52%  "

$ws.Cells.Item(42, 5).Value = "  -3.95%  "

$ws.Cells.Item(43, 4).Value = "1.949.53"
$ws.Cells.Item(43, 5).Value = "  -0.75%  "

$ws.Cells.Item(44, 5).Value = "  -3.72%  "

$ws.Cells.Item(45, 5).Value = "  +0.83%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "9.39"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.63%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "16.16"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -8.79%  "

$ws.Cells.Item(48, 5).Value = "  -4.38%  "

$ws.Cells.Item(49, 4).Value = "2.441.51"
$ws.Cells.Item(49, 5).Value = "  -4.87%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "70.82"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.30%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "89.05"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -4.53%  "
